$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: translate column titles to snake_case English names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Normalize "de" -> "De" in a handful of place names
$ws.Range("A2").Value = "Ciudad De México"
$ws.Range("A4").Value = "Estado De México"
$ws.Range("B4").Value = "Naucalpan De Juárez"
$ws.Range("A7").Value = "Guanajuato"
$ws.Range("B12").Value = "Progreso De Obregón"
$ws.Range("B16").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B22").Value = "Villa De Zaachila"

# Remove the trailing footer/metadata rows (30-34)
$ws.Range("A30:D34").EntireRow.Delete()
